$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.551.39"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'2.043.43"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'245.33"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'56.53"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'63.19"
$ws.Range("E9").Value = "  +6.52%  "
$ws.Range("D10").Value = "'0.369"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "'0.0748"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").Value = "'0.106"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").Value = "'0.905"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "'14.28"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "'2.347.35"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "'5.40"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").Value = "'2.068.13"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "'17.74"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "'36.455.34"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "'71.65"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "'0.0₃0856"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").Value = "'236.70"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'5.20"
$ws.Range("E23").Value = "  -5.01%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +3.82%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  -6.85%  "
$ws.Range("D28").Value = "'163.91"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").Value = "'19.93"
$ws.Range("E29").Value = "  -4.96%  "
$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "'1.22"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").Value = "'4.98"
$ws.Range("E32").Value = "  -6.05%  "
$ws.Range("D33").Value = "'0.0599"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").Value = "'4.40"
$ws.Range("E34").Value = "  -6.60%  "
$ws.Range("D35").Value = "'0.0878"
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").Value = "'2.20"
$ws.Range("E38").Value = "  -7.41%  "
$ws.Range("D39").Value = "'5.11"
$ws.Range("E39").Value = "  +3.54%  "
$ws.Range("D40").Value = "'1.22"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("D41").Value = "'0.0216"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "'2.88"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").Value = "'1.10"
$ws.Range("E43").Value = "  -4.72%  "
$ws.Range("D44").Value = "'93.62"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "'0.0905"
$ws.Range("E45").Value = "  -5.28%  "
$ws.Range("D46").Value = "'15.92"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "'1.373.90"
$ws.Range("E47").Value = "  +2.55%  "
$ws.Range("D48").Value = "'7.43"
$ws.Range("E48").Value = "  +6.13%  "
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("D50").Value = "'2.26"
$ws.Range("E50").Value = "  -6.22%  "
$ws.Range("D51").Value = "'45.97"
$ws.Range("E51").Value = "  +0.47%  "
